$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update footer timestamp text
$ws.Range("A1").Value = "Datos actualizados a 8 de Septiembre de 2020 a las 19:13"

# Country row A128/A129: Gambia moves above Siria (new Gambia totals overtake old Siria total),
# Siria keeps its previous (unchanged) stats one row lower.
$ws.Range("A128").Value = "Gambia"
$ws.Range("A129").Value = "Siria"

# Refreshed daily COVID figures per country row
# Row 4
$ws.Range("B4").Value = 6493058
$ws.Range("C4").Value = 7483
$ws.Range("D4").Value = 3760821
$ws.Range("E4").Value = 2538564
$ws.Range("G4").Value = 139
$ws.Range("H4").Value = 193673

# Row 5
$ws.Range("B5").Value = 4338267
$ws.Range("C5").Value = 60683
$ws.Range("D5").Value = 3377530
$ws.Range("E5").Value = 887235
$ws.Range("G5").Value = 686
$ws.Range("H5").Value = 73502

# Row 6
$ws.Range("B6").Value = 4150311
$ws.Range("C6").Value = 2517
$ws.Range("E6").Value = 667663
$ws.Range("G6").Value = 83
$ws.Range("H6").Value = 127084

# Row 14
$ws.Range("B14").Value = 425541
$ws.Range("C14").Value = 1267
$ws.Range("D14").Value = 397730
$ws.Range("E14").Value = 16129
$ws.Range("G14").Value = 30
$ws.Range("H14").Value = 11682

# Row 16
$ws.Range("B16").Value = 352520
$ws.Range("C16").Value = 2420
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = 41584

# Row 21
$ws.Range("B21").Value = 283270
$ws.Range("C21").Value = 1761
$ws.Range("D21").Value = 253245
$ws.Range("E21").Value = 23243
$ws.Range("G21").Value = 52
$ws.Range("H21").Value = 6782

# Row 59
$ws.Range("B59").Value = 46938
$ws.Range("C59").Value = 285
$ws.Range("D59").Value = 33183
$ws.Range("E59").Value = 12184
$ws.Range("G59").Value = 9
$ws.Range("H59").Value = 1571

# Row 71
$ws.Range("B71").Value = 30080
$ws.Range("C71").Value = 306
$ws.Range("E71").Value = 4938
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 1778

# Row 72
$ws.Range("B72").Value = 29368
$ws.Range("C72").Value = 652
$ws.Range("D72").Value = 20139
$ws.Range("E72").Value = 8789
$ws.Range("G72").Value = 3
$ws.Range("H72").Value = 440

# Row 79
$ws.Range("B79").Value = 21324
$ws.Range("C79").Value = 498
$ws.Range("D79").Value = 6722
$ws.Range("E79").Value = 14395
$ws.Range("G79").Value = 7
$ws.Range("H79").Value = 207

# Row 101
$ws.Range("B101").Value = 8741
$ws.Range("C101").Value = 74
$ws.Range("D101").Value = 6157
$ws.Range("E101").Value = 2555

# Row 119
$ws.Range("B119").Value = 4647
$ws.Range("C119").Value = 90
$ws.Range("D119").Value = 2715
$ws.Range("E119").Value = 1904
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 28

# Row 128
$ws.Range("B128").Value = 3275
$ws.Range("C128").Value = 78
$ws.Range("D128").Value = 1424
$ws.Range("E128").Value = 1752
$ws.Range("H128").Value = 99

# Row 129
$ws.Range("B129").Value = 3229
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 744
$ws.Range("E129").Value = 2348
$ws.Range("H129").Value = 137
